# Backlog Doc.xlsx - small changes to sprint log
#
# Sprint 1 Backlog: the two "QT library" rows (13 & 14) move from
# "On going"/"?" to "Not Completed" with no finish date yet.
#
# Sprint 2 Backlog: a few week / total hours get bumped up, and the
# sprint summary row's week totals + grand total are corrected.

$wb = $excel.ActiveWorkbook

# --- Sprint 1 Backlog ---
$ws1 = $wb.Worksheets.Item("Sprint 1 Backlog")

# Row 13: "Start implementing QT library to add a GUI to the project"
$ws1.Cells.Item(13, 11).Value = $null          # K13 was "?"
$ws1.Cells.Item(13, 12).Value = "Not Completed" # L13 was "On going"

# Row 14: "Start Creating grahpical interface by using QT"
$ws1.Cells.Item(14, 11).Value = $null          # K14 was "?"
$ws1.Cells.Item(14, 12).Value = "Not Completed" # L14 was "On going"

# --- Sprint 2 Backlog ---
$ws2 = $wb.Worksheets.Item("Sprint 2 Backlog")

# Row 9: "Start the functionality of the board/level..."
$ws2.Cells.Item(9, 7).Value = 15   # G9: 11 -> 15
$ws2.Cells.Item(9, 10).Value = 15  # J9: 11 -> 15

# Row 10: "Start the functionality of the snake..."
$ws2.Cells.Item(10, 9).Value = 5   # I10: 2 -> 5
$ws2.Cells.Item(10, 10).Value = 15 # J10: 12 -> 15

# Row 12: "Start the functionality of the menu..."
$ws2.Cells.Item(12, 8).Value = 10  # H12: 9 -> 10
$ws2.Cells.Item(12, 10).Value = 15 # J12: 14 -> 15

# Row 19: Completed Sprint summary row
$ws2.Cells.Item(19, 7).Value = 15  # G19: 30 -> 15
$ws2.Cells.Item(19, 8).Value = 30  # H19: (blank) -> 30
$ws2.Cells.Item(19, 9).Value = 20  # I19: (blank) -> 20
$ws2.Cells.Item(19, 10).Value = 65 # J19: 57 -> 65
